# PowerShell (Excel COM-interop) script implementing the
# "25/01/2017 Update after GW call" edit to Gemalto-CloudGate-Issues_O.xlsx
#
# Summary of changes on the "Issue Tracking" sheet:
#  - View zoomed from 90% to 85%, active selection moved to F6
#  - Row 6 (issue #4): now Closed, hidden by the filter, F6 gets a new comment
#  - Row 7 (issue #5): E7 resolution note gets a new bullet point appended
#  - Row 11 (issue #9): now Closed, hidden by the filter, F11 gets a new comment
#  - Row 12 (issue #10): F12 comment replaced with a new 01/24 comment
#  - Row 13 (issue #11): F13 comment replaced with the same new 01/24 comment
#  - AutoFilter's remembered checked values updated to match the now-visible rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue Tracking")

# --- Row 6 (issue #4) ---------------------------------------------------
$ws.Range("F6").Value = '•01/24 Brandon sees no issue in testing with 2.71.'
$ws.Range("G6").Value = "Closed"

# --- Row 7 (issue #5) ----------------------------------------------------
$ws.Range("E7").Value = '•This is the same on both the current CG 3G (Gobi) and the CG LTE. There was no change compared to those devices.
•We will try to get to a solution for the release following the 2.70.0. Some testing will be needed from GetWireless and Option.
•Engineering build will be provided in January for testing and if OK, then a release will happen the same month.
•01/06 Engineering build provided to GetWireless. GetWireless to test.
•01/11 It will be great if we can get a result by the end of the week so that we can add it on the release candidate.
•01/18 Jimmy is looking into this one.
•01/25 There seems to be a problem with the modem that stays on SIM not ready and that is why is showing the problem after a radio modem reboot. We need to look ino the modem log and either fix it ourselves or get a fix from Gemalto.'

# --- Row 11 (issue #9) ----------------------------------------------------
$ws.Range("F11").Value = '•01/24 The modem does not auto detect the SIM and load the UMTS firmware. The modem default must be UMTS causing this behavior. The gobi default was Verizon.'
$ws.Range("G11").Value = "Closed"

# --- Row 12 (issue #10) ---------------------------------------------------
$ws.Range("F12").Value = '•01/24 Brandon saw no problem with the firewall and LuvitRED opening the same port.'

# --- Row 13 (issue #11) ---------------------------------------------------
$ws.Range("F13").Value = '•01/24 Brandon saw no problem with the firewall and LuvitRED opening the same port.'

# --- Hide rows 6 and 11 (now filtered out as "Closed") --------------------
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(11).Hidden = $true

# --- Update the AutoFilter's remembered checked values ---------------------
$critValues = @("Status", "Under Investigation (OP)", "Under Test (GW)", "")
$ws.Range("G1:G14").AutoFilter(1, $critValues, 7)

# --- View changes: zoom to 85% and move selection to F6 --------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("F6").Select()
